$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New titration/CRM check logged for the tank opened 2022-04-27 (row 83),
# continuing the D column's "% off" shared formula pattern from the rows above.
$ws.Range("A83").Value = 20220427
$ws.Range("B83").Value = 2215.5293579999998
$ws.Range("C83").Value = 2224.4699999999998
$ws.Range("D83").Formula = "=100*(B83-C83)/C83"
$ws.Range("E83").Value = 180
$ws.Range("F83").Value = "CRM OPENED 20220427"

# Mirror the post-edit view state: user had scrolled down and landed on C88.
$ws.Range("C88").Select()
